$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update cell values per the diff:
# B1: "Fun" -> "asd"
# C1: 1 -> 123
# B2: 23 -> "asd"
$ws.Range("B1").Value = "asd"
$ws.Range("C1").Value = 123
$ws.Range("B2").Value = "asd"

# Minimize the workbook window
$excel.ActiveWindow.WindowState = -4140
